# Auto-generated edit script applying numeric corrections to Titan_Profits data
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (sheet1..sheet8).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # row 15
$ws.Range("H15").Value = 99308.55
$ws.Range("I15").Value = 99308.55
$ws.Range("K15").Value = 297925.65
$ws.Range("M15").Value = -297756.65

$ws = $wb.Worksheets.Item(1)  # row 137
$ws.Range("H137").Value = 32259130
$ws.Range("I137").Value = 33334368
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 100003104
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -100000554
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item(2)  # row 2
$ws.Range("H2").Value = 127115.25
$ws.Range("I2").Value = 145103.14
$ws.Range("K2").Value = 145103.14
$ws.Range("M2").Value = -144990.14

$ws = $wb.Worksheets.Item(2)  # row 32
$ws.Range("H32").Value = 22474.213
$ws.Range("I32").Value = 5121.281
$ws.Range("J32").Value = 269753.5
$ws.Range("K32").Value = 5121.281
$ws.Range("L32").Value = 269753.5
$ws.Range("M32").Value = -4834.281
$ws.Range("N32").Value = -270327.5

$ws = $wb.Worksheets.Item(2)  # row 61
$ws.Range("H61").Value = 1507.7428
$ws.Range("I61").Value = 1255.3438
$ws.Range("K61").Value = 1255.3438
$ws.Range("M61").Value = -1043.3438

$ws = $wb.Worksheets.Item(2)  # row 74
$ws.Range("H74").Value = 3432.22
$ws.Range("I74").Value = 957.46155
$ws.Range("J74").Value = 12206.363
$ws.Range("K74").Value = 957.46155
$ws.Range("L74").Value = 12206.363
$ws.Range("M74").Value = -83.46154999999999
$ws.Range("N74").Value = -13954.363

$ws = $wb.Worksheets.Item(2)  # row 77
$ws.Range("H77").Value = 3432.22
$ws.Range("I77").Value = 957.46155
$ws.Range("J77").Value = 12206.363
$ws.Range("K77").Value = 4787.30775
$ws.Range("L77").Value = 61031.815
$ws.Range("M77").Value = -419.3077499999999
$ws.Range("N77").Value = -69767.815

$ws = $wb.Worksheets.Item(2)  # row 116
$ws.Range("H116").Value = 127115.25
$ws.Range("I116").Value = 145103.14
$ws.Range("K116").Value = 145103.14
$ws.Range("M116").Value = -142809.14

$ws = $wb.Worksheets.Item(2)  # row 132
$ws.Range("H132").Value = 1854.6774
$ws.Range("I132").Value = 1647.537
$ws.Range("J132").Value = 3252.875
$ws.Range("K132").Value = 4942.611
$ws.Range("L132").Value = 9758.625
$ws.Range("M132").Value = -2412.611
$ws.Range("N132").Value = -14818.625

$ws = $wb.Worksheets.Item(2)  # row 136
$ws.Range("H136").Value = 1507.7428
$ws.Range("I136").Value = 1255.3438
$ws.Range("K136").Value = 3766.0314
$ws.Range("M136").Value = -1216.0314

$ws = $wb.Worksheets.Item(3)  # row 3
$ws.Range("H3").Value = 127115.25
$ws.Range("I3").Value = 145103.14
$ws.Range("K3").Value = 145103.14
$ws.Range("M3").Value = -144989.14

$ws = $wb.Worksheets.Item(3)  # row 105
$ws.Range("H105").Value = 3672.5454
$ws.Range("I105").Value = 3581.5454
$ws.Range("J105").Value = 3854.5454
$ws.Range("K105").Value = 3581.5454
$ws.Range("L105").Value = 3854.5454
$ws.Range("M105").Value = -1834.5454
$ws.Range("N105").Value = -7348.5454

$ws = $wb.Worksheets.Item(3)  # row 119
$ws.Range("H119").Value = 29500
$ws.Range("J119").Value = 29500
$ws.Range("L119").Value = 29500
$ws.Range("N119").Value = -39176

$ws = $wb.Worksheets.Item(3)  # row 134
$ws.Range("H134").Value = 27780666
$ws.Range("I134").Value = 34484732
$ws.Range("J134").Value = 6677.857
$ws.Range("K134").Value = 103454196
$ws.Range("L134").Value = 20033.571
$ws.Range("M134").Value = -103451661
$ws.Range("N134").Value = -25103.571

$ws = $wb.Worksheets.Item(4)  # row 31
$ws.Range("H31").Value = 3723.8823
$ws.Range("I31").Value = 1808.0869
$ws.Range("J31").Value = 5297.5713
$ws.Range("K31").Value = 1808.0869
$ws.Range("L31").Value = 5297.5713
$ws.Range("M31").Value = -1513.0869
$ws.Range("N31").Value = -5887.5713

$ws = $wb.Worksheets.Item(4)  # row 34
$ws.Range("H34").Value = 3723.8823
$ws.Range("I34").Value = 1808.0869
$ws.Range("J34").Value = 5297.5713
$ws.Range("K34").Value = 1808.0869
$ws.Range("L34").Value = 5297.5713
$ws.Range("M34").Value = -1606.0869
$ws.Range("N34").Value = -5701.5713

$ws = $wb.Worksheets.Item(4)  # row 58
$ws.Range("H58").Value = 27779482
$ws.Range("I58").Value = 37038100
$ws.Range("J58").Value = 3623.2222
$ws.Range("K58").Value = 37038100
$ws.Range("L58").Value = 3623.2222
$ws.Range("M58").Value = -37037897
$ws.Range("N58").Value = -4029.2222

$ws = $wb.Worksheets.Item(4)  # row 132
$ws.Range("H132").Value = 3473638
$ws.Range("I132").Value = 4505426
$ws.Range("J132").Value = 3077.3635
$ws.Range("K132").Value = 13516278
$ws.Range("L132").Value = 9232.0905
$ws.Range("M132").Value = -13513748
$ws.Range("N132").Value = -14292.0905

$ws = $wb.Worksheets.Item(4)  # row 136
$ws.Range("H136").Value = 27779482
$ws.Range("I136").Value = 37038100
$ws.Range("J136").Value = 3623.2222
$ws.Range("K136").Value = 111114300
$ws.Range("L136").Value = 10869.6666
$ws.Range("M136").Value = -111111750
$ws.Range("N136").Value = -15969.6666

$ws = $wb.Worksheets.Item(5)  # row 33
$ws.Range("H33").Value = 83.833336
$ws.Range("I33").Value = 64.59999999999999
$ws.Range("J33").Value = 180
$ws.Range("K33").Value = 387.6
$ws.Range("L33").Value = 1080
$ws.Range("M33").Value = -104.6
$ws.Range("N33").Value = -1646

$ws = $wb.Worksheets.Item(5)  # row 34
$ws.Range("H34").Value = 3726.8
$ws.Range("J34").Value = 3671.4285
$ws.Range("L34").Value = 11014.2855
$ws.Range("N34").Value = -11182.2855

$ws = $wb.Worksheets.Item(5)  # row 39
$ws.Range("H39").Value = 9220
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 9220
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 27660
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -28248

$ws = $wb.Worksheets.Item(5)  # row 55
$ws.Range("H55").Value = 4208.3335
$ws.Range("J55").Value = 4545.4546
$ws.Range("L55").Value = 13636.3638
$ws.Range("N55").Value = -13990.3638

$ws = $wb.Worksheets.Item(5)  # row 82
$ws.Range("H82").Value = 3985.7144
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3985.7144
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 11957.1432
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -12769.1432

$ws = $wb.Worksheets.Item(5)  # row 85
$ws.Range("H85").Value = 3985.7144
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3985.7144
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 11957.1432
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -14765.1432

$ws = $wb.Worksheets.Item(5)  # row 113
$ws.Range("H113").Value = 12987620
$ws.Range("J113").Value = 12987620
$ws.Range("L113").Value = 38962860
$ws.Range("N113").Value = -38967200

$ws = $wb.Worksheets.Item(5)  # row 127
$ws.Range("H127").Value = 1354.1177
$ws.Range("J127").Value = 1401.25
$ws.Range("L127").Value = 4203.75
$ws.Range("N127").Value = -14123.75

$ws = $wb.Worksheets.Item(5)  # row 131
$ws.Range("H131").Value = 9525671
$ws.Range("I131").Value = 566.6667
$ws.Range("J131").Value = 10418650
$ws.Range("K131").Value = 1700.0001
$ws.Range("L131").Value = 31255950
$ws.Range("M131").Value = 3339.9999
$ws.Range("N131").Value = -31266030

$ws = $wb.Worksheets.Item(5)  # row 134
$ws.Range("H134").Value = 6254.9414
$ws.Range("I134").Value = 4027.8333
$ws.Range("J134").Value = 11600
$ws.Range("K134").Value = 12083.4999
$ws.Range("L134").Value = 34800
$ws.Range("M134").Value = -7013.499899999999
$ws.Range("N134").Value = -44940

$ws = $wb.Worksheets.Item(6)  # row 102
$ws.Range("H102").Value = 3450.6758
$ws.Range("I102").Value = 2246.7036
$ws.Range("J102").Value = 6701.4
$ws.Range("K102").Value = 2246.7036
$ws.Range("L102").Value = 6701.4
$ws.Range("M102").Value = -624.7035999999998
$ws.Range("N102").Value = -9945.4

$ws = $wb.Worksheets.Item(6)  # row 113
$ws.Range("H113").Value = 1607.2858
$ws.Range("I113").Value = 1465.6364
$ws.Range("J113").Value = 2126.6667
$ws.Range("K113").Value = 1465.6364
$ws.Range("L113").Value = 2126.6667
$ws.Range("M113").Value = 704.3635999999999
$ws.Range("N113").Value = -6466.6667

$ws = $wb.Worksheets.Item(6)  # row 132
$ws.Range("H132").Value = 4013.7585
$ws.Range("I132").Value = 3567.4546
$ws.Range("J132").Value = 5416.4287
$ws.Range("K132").Value = 10702.3638
$ws.Range("L132").Value = 16249.2861
$ws.Range("M132").Value = -8172.363799999999
$ws.Range("N132").Value = -21309.2861

$ws = $wb.Worksheets.Item(6)  # row 138
$ws.Range("H138").Value = 71803.22
$ws.Range("J138").Value = 71803.22
$ws.Range("L138").Value = 71803.22
$ws.Range("N138").Value = -82083.22

$ws = $wb.Worksheets.Item(7)  # row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws = $wb.Worksheets.Item(7)  # row 136
$ws.Range("H136").Value = 8037.4443
$ws.Range("I136").Value = 5659.778
$ws.Range("J136").Value = 10415.111
$ws.Range("K136").Value = 16979.334
$ws.Range("L136").Value = 31245.333
$ws.Range("M136").Value = -14429.334
$ws.Range("N136").Value = -36345.333

$ws = $wb.Worksheets.Item(8)  # row 132
$ws.Range("H132").Value = 3089.6924
$ws.Range("I132").Value = 2927.147
$ws.Range("J132").Value = 4195
$ws.Range("K132").Value = 8781.440999999999
$ws.Range("L132").Value = 12585
$ws.Range("M132").Value = -6251.440999999999
$ws.Range("N132").Value = -17645

$ws = $wb.Worksheets.Item(8)  # row 136
$ws.Range("H136").Value = 3510.4614
$ws.Range("I136").Value = 1511.0435
$ws.Range("J136").Value = 6384.625
$ws.Range("K136").Value = 4533.1305
$ws.Range("L136").Value = 19153.875
$ws.Range("M136").Value = -1983.1305
$ws.Range("N136").Value = -24253.875
